$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: min
$ws.Range("A22").Value = "min"
$ws.Range("A22").Font.Bold = $true
$ws.Range("B22").Formula = "=MIN(B3:B20)"
$ws.Range("C22").Formula = "=MIN(C3:C20)"

# Row 23: max
$ws.Range("A23").Value = "max"
$ws.Range("A23").Font.Bold = $true
$ws.Range("B23").Formula = "=MAX(B3:B20)"
$ws.Range("C23").Formula = "=MAX(C3:C20)"

# Row 24: sum
$ws.Range("A24").Value = "sum"
$ws.Range("A24").Font.Bold = $true
$ws.Range("B24").Formula = "=SUM(B3:B20)"
$ws.Range("C24").Formula = "=SUM(C3:C20)"

# Update selection to match the new active cell/selection block
$null = $ws.Range("B24:C24").Select()
